# "Add files via upload" — the uploaded workbook now has the
# "Approved/Rejected" column (I) filled in with "Approved" for each of
# the four visible TestScenario summary rows (2, 5, 7, 10); the hidden
# detail rows in between are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "Approved"
$ws.Range("I5").Value = "Approved"
$ws.Range("I7").Value = "Approved"
$ws.Range("I10").Value = "Approved"

# Reflect where the user ended up looking: scrolled right so column G is
# the first visible column, with I2:I10 selected (active cell I2) —
# i.e. the column they had just finished filling in.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I2:I10").Select() | Out-Null
